# Update LR-pair TPM-derived statistics (ligand/receptor expression
# values, specificities and edge weights) with newly recomputed TPM
# values, per "update scripts wuth new tpm".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("G2").Value = 13.32779766666667
$ws.Range("H2").Value = 39.983393
$ws.Range("I2").Value = 0.1697233513642653
$ws.Range("J2").Value = 0.1697233513642653
$ws.Range("M2").Value = 0.5282606666666667
$ws.Range("N2").Value = 1.584782
$ws.Range("O2").Value = 0.04843445786516468
$ws.Range("P2").Value = 0.04843445786516468
$ws.Range("Q2").Value = 7.040551280591779
$ws.Range("R2").Value = 63.364961525326
$ws.Range("S2").Value = 0.008220458510387049
$ws.Range("T2").Value = 0.008220458510387049
# Row 3: ECs -> FAPs
$ws.Range("G3").Value = 13.32779766666667
$ws.Range("H3").Value = 39.983393
$ws.Range("I3").Value = 0.1697233513642653
$ws.Range("J3").Value = 0.1697233513642653
$ws.Range("M3").Value = 8.934372
$ws.Range("N3").Value = 26.803116
$ws.Range("O3").Value = 0.8191627571218761
$ws.Range("P3").Value = 0.8191627571218761
$ws.Range("Q3").Value = 119.075502294732
$ws.Range("R3").Value = 1071.679520652588
$ws.Range("S3").Value = 0.1390310484515165
$ws.Range("T3").Value = 0.1390310484515165
# Row 4: ECs -> MuSCs
$ws.Range("G4").Value = 13.32779766666667
$ws.Range("H4").Value = 39.983393
$ws.Range("I4").Value = 0.1697233513642653
$ws.Range("J4").Value = 0.1697233513642653
$ws.Range("M4").Value = 1.444079
$ws.Range("N4").Value = 4.332237
$ws.Range("O4").Value = 0.1324027850129592
$ws.Range("P4").Value = 0.1324027850129592
$ws.Range("Q4").Value = 19.24639272668233
$ws.Range("R4").Value = 173.217534540141
$ws.Range("S4").Value = 0.02247184440236175
$ws.Range("T4").Value = 0.02247184440236175
# Row 5: FAPs -> ECs
$ws.Range("G5").Value = 43.30706799999999
$ws.Range("H5").Value = 129.921204
$ws.Range("I5").Value = 0.5514955210569645
$ws.Range("J5").Value = 0.5514955210569645
$ws.Range("M5").Value = 0.5282606666666667
$ws.Range("N5").Value = 1.584782
$ws.Range("O5").Value = 0.04843445786516468
$ws.Range("P5").Value = 0.04843445786516468
$ws.Range("Q5").Value = 22.87742061305866
$ws.Range("R5").Value = 205.896785517528
$ws.Range("S5").Value = 0.02671138657746059
$ws.Range("T5").Value = 0.02671138657746059
# Row 6: FAPs -> FAPs
$ws.Range("G6").Value = 43.30706799999999
$ws.Range("H6").Value = 129.921204
$ws.Range("I6").Value = 0.5514955210569645
$ws.Range("J6").Value = 0.5514955210569645
$ws.Range("M6").Value = 8.934372
$ws.Range("N6").Value = 26.803116
$ws.Range("O6").Value = 0.8191627571218761
$ws.Range("P6").Value = 0.8191627571218761
$ws.Range("Q6").Value = 386.9214557412959
$ws.Range("R6").Value = 3482.293101671663
$ws.Range("S6").Value = 0.4517645915693888
$ws.Range("T6").Value = 0.4517645915693888
# Row 7: FAPs -> MuSCs
$ws.Range("G7").Value = 43.30706799999999
$ws.Range("H7").Value = 129.921204
$ws.Range("I7").Value = 0.5514955210569645
$ws.Range("J7").Value = 0.5514955210569645
$ws.Range("M7").Value = 1.444079
$ws.Range("N7").Value = 4.332237
$ws.Range("O7").Value = 0.1324027850129592
$ws.Range("P7").Value = 0.1324027850129592
$ws.Range("Q7").Value = 62.538827450372
$ws.Range("R7").Value = 562.849447053348
$ws.Range("S7").Value = 0.07301954291011517
$ws.Range("T7").Value = 0.07301954291011517
# Row 8: MuSCs -> ECs
$ws.Range("G8").Value = 21.891734
$ws.Range("H8").Value = 65.67520200000001
$ws.Range("I8").Value = 0.27878112757877
$ws.Range("J8").Value = 0.27878112757877
$ws.Range("M8").Value = 0.5282606666666667
$ws.Range("N8").Value = 1.584782
$ws.Range("O8").Value = 0.04843445786516468
$ws.Range("P8").Value = 0.04843445786516468
$ws.Range("Q8").Value = 11.56454199732934
$ws.Range("R8").Value = 104.080877975964
$ws.Range("S8").Value = 0.01350261277731704
$ws.Range("T8").Value = 0.01350261277731704
# Row 9: MuSCs -> FAPs
$ws.Range("G9").Value = 21.891734
$ws.Range("H9").Value = 65.67520200000001
$ws.Range("I9").Value = 0.27878112757877
$ws.Range("J9").Value = 0.27878112757877
$ws.Range("M9").Value = 8.934372
$ws.Range("N9").Value = 26.803116
$ws.Range("O9").Value = 0.8191627571218761
$ws.Range("P9").Value = 0.8191627571218761
$ws.Range("Q9").Value = 195.588895281048
$ws.Range("R9").Value = 1760.300057529432
$ws.Range("S9").Value = 0.2283671171009708
$ws.Range("T9").Value = 0.2283671171009708
# Row 10: MuSCs -> MuSCs
$ws.Range("G10").Value = 21.891734
$ws.Range("H10").Value = 65.67520200000001
$ws.Range("I10").Value = 0.27878112757877
$ws.Range("J10").Value = 0.27878112757877
$ws.Range("M10").Value = 1.444079
$ws.Range("N10").Value = 4.332237
$ws.Range("O10").Value = 0.1324027850129592
$ws.Range("P10").Value = 0.1324027850129592
$ws.Range("Q10").Value = 31.61339334298601
$ws.Range("R10").Value = 284.5205400868741
$ws.Range("S10").Value = 0.03691139770048223
$ws.Range("T10").Value = 0.03691139770048223